# Adding the changes we made on may 9th
# 1) Shift existing data rows (2..21) down by 4 rows (to 6..25) to make room
#    for 4 new rows at the top of the data, reading with Value2 (Value's
#    getter is not reliable in this host) and writing with Value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 21; $r -ge 2; $r--) {
    $dest = $r + 4
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# 2) Fill in the 4 new rows inserted at the top (rows 2-5)
$ws.Cells.Item(2, 1).Value = -1.238285183906555
$ws.Cells.Item(2, 2).Value = 0.4941155612468719
$ws.Cells.Item(2, 3).Value = -1.145785808563232

$ws.Cells.Item(3, 1).Value = 0.3425590991973889
$ws.Cells.Item(3, 2).Value = 0.401696681976318
$ws.Cells.Item(3, 3).Value = -1.307540893554688

$ws.Cells.Item(4, 1).Value = 1.148096084594726
$ws.Cells.Item(4, 2).Value = 0.2248815298080447
$ws.Cells.Item(4, 3).Value = -1.84254863858223

$ws.Cells.Item(5, 1).Value = 0.3787193298339841
$ws.Cells.Item(5, 2).Value = 0.750096321105957
$ws.Cells.Item(5, 3).Value = -2.375997304916381

# 3) Append 6 new rows of data at the bottom (rows 26-31)
$ws.Cells.Item(26, 1).Value = -0.04362952709197906
$ws.Cells.Item(26, 2).Value = -0.3616583049297337
$ws.Cells.Item(26, 3).Value = -1.960778713226318

$ws.Cells.Item(27, 1).Value = 0.191988468170166
$ws.Cells.Item(27, 2).Value = -0.4578718543052673
$ws.Cells.Item(27, 3).Value = -1.877070605754853

$ws.Cells.Item(28, 1).Value = 0.1683353185653686
$ws.Cells.Item(28, 2).Value = -0.3762182295322418
$ws.Cells.Item(28, 3).Value = -1.764573842287064

$ws.Cells.Item(29, 1).Value = 0.1497325897216795
$ws.Cells.Item(29, 2).Value = -0.3715704679489131
$ws.Cells.Item(29, 3).Value = -1.762347698211669

$ws.Cells.Item(30, 1).Value = 0.0995370149612424
$ws.Cells.Item(30, 2).Value = -0.2668604403734202
$ws.Cells.Item(30, 3).Value = -1.60942207276821

$ws.Cells.Item(31, 1).Value = 0.02663779258727977
$ws.Cells.Item(31, 2).Value = -0.2185956239700315
$ws.Cells.Item(31, 3).Value = -1.59472194314003
